$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 becomes a "separator" row (like row 9): copy format from row 9 ---
$ws.Range("A9:E9").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# --- Rows 12-13: same banding as rows 3-7 (thin top+bottom border, small font) ---
$ws.Range("A3:E3").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)

# --- Rows 14-17: same banding as row 8/10 (no border, small font) ---
$ws.Range("A8:E8").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
# Rows 15-17 have no cell at all in column A, so only paste B:E
$ws.Range("B8:E8").Copy()
$ws.Range("B15:E15").PasteSpecial(-4122)
$ws.Range("B8:E8").Copy()
$ws.Range("B16:E16").PasteSpecial(-4122)
$ws.Range("B8:E8").Copy()
$ws.Range("B17:E17").PasteSpecial(-4122)

# --- Row 12 ---
$ws.Range("A12").Value = "SCRIPT/T01P01A/us0106.ssb"
$ws.Range("B12").Value = 225
$ws.Range("C12").Value = " You two are great! ♪"
$ws.Range("D12").Value = " Вы двое просто великолепны! ♪"
$ws.Range("E12").Value = " Âú äâïå ðñïòóï âåìéëïìåðîú! ♪"

# --- Row 13 ---
$ws.Range("A13").Value = "SCRIPT/P01P04A/us3101.ssb"
$ws.Range("B13").Value = 206
$ws.Range("C13").Value = " Let\'s find some treasure... ♪"
$ws.Range("D13").Value = " Пора найти сокровища... ♪"
$ws.Range("E13").Value = " Ðïñà îàêóé òïëñïâéþà... ♪"

# --- Row 14 ---
$ws.Range("A14").Value = "SCRIPT/D73P11A/us3121.ssb"
$ws.Range("B14").Value = 187
$ws.Range("C14").Value = " [CS:N]Shaymin[CR] are so small\nand cute. ♪"
$ws.Range("D14").Value = " [CS:N]Шеймины[CR] такие маленькие и\nмиленькие. ♪"
$ws.Range("E14").Value = " [CS:N]Šåêíéîú[CR] óàëéå íàìåîûëéå é\níéìåîûëéå. ♪"

# --- Row 15 ---
$ws.Range("B15").Value = 165
$ws.Range("C15").Value = " It\'s a so-called secret, so it \nmust be hidden…"
$ws.Range("D15").Value = " Её называют секретной, значит\nона скрыта..."
$ws.Range("E15").Value = " Åæ îàèúâàýó òåëñåóîïê, èîàœéó\nïîà òëñúóà..."

# --- Row 16 ---
$ws.Range("B16").Value = 168
$ws.Range("C16").Value = " It\'s hard to find, huh?\nHa ha ha. ♪"
$ws.Range("D16").Value = " Её непросто найти, а?\nХа-ха-ха. ♪"
$ws.Range("E16").Value = " Åæ îåðñïòóï îàêóé, à?\nÖà-öà-öà. ♪"

# --- Row 17 ---
$ws.Range("B17").Value = 146
$ws.Range("C17").Value = " This place is so relaxing. ♪"
$ws.Range("D17").Value = " Здесь так хорошо. ♪"
$ws.Range("E17").Value = " Èäåòû óàë öïñïšï. ♪"

# --- Row heights (row 11 already carries its natural 43.2 height from before) ---
$ws.Range("A12:E12").RowHeight = 43.2
$ws.Range("A13:E13").RowHeight = 43.2
$ws.Range("A14:E14").RowHeight = 43.2
$ws.Range("A15:E15").RowHeight = 21.6
$ws.Range("A16:E16").RowHeight = 25.2

# --- View state ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D18").Select()
